{"js": "// \"Deje como estaba el plan de iteracion\"\n// 1) \"Metricas\" -> \"M\u00e9tricas\" (remove spell-check proofErr markers, add accent)\n// 2) \"Herramientas a utilizar\" -> split into \"Describir h\" + \"erramientas a utilizar\"\n//    (visible text becomes \"Describir herramientas a utilizar\")\n// 3) \"Descripcion\" -> \"Descripci\u00f3n\" (remove spell-check proofErr markers, add accent)\n// 4) \"Responsables\" -> \"responsables\" (remove grammar-check proofErr markers, lowercase)\n\nconst PKG_OPEN = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>';\nconst PKG_CLOSE = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\n// Replace the whole paragraph that contains `searchText` with a freshly built\n// paragraph made of `runsXml` (an array of already-serialized <w:r>...</w:r> strings).\n// Using the *whole paragraph* range (not just the narrow search hit) ensures any\n// sibling <w:proofErr/> spell/grammar markers inside that paragraph are dropped too.\nasync function replaceParagraphRuns(context, searchText, runsXml) {\n  const body = context.document.body;\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly one match for \"${searchText}\", found ${results.items.length}`);\n  }\n\n  const hit = results.items[0];\n  const para = hit.paragraphs.getFirst();\n  const wholeParaRange = para.getRange(\"Whole\");\n\n  const ooxml = PKG_OPEN + \"<w:p>\" + runsXml.join(\"\") + \"</w:p>\" + PKG_CLOSE;\n  wholeParaRange.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Metricas -> M\u00e9tricas\nawait replaceParagraphRuns(context, \"Metricas\", [\n  \"<w:r><w:t>M\u00e9tricas</w:t></w:r>\",\n  '<w:r><w:t xml:space=\"preserve\"> basadas en casos de uso</w:t></w:r>'\n]);\n\n// 2) Herramientas a utilizar -> Describir h | erramientas a utilizar\nawait replaceParagraphRuns(context, \"Herramientas a utilizar\", [\n  \"<w:r><w:t>Describir h</w:t></w:r>\",\n  \"<w:r><w:t>erramientas a utilizar</w:t></w:r>\"\n]);\n\n// 3) Descripcion -> Descripci\u00f3n\nawait replaceParagraphRuns(context, \"Descripcion\", [\n  \"<w:r><w:t>Descripci\u00f3n</w:t></w:r>\",\n  '<w:r><w:t xml:space=\"preserve\"> general</w:t></w:r>'\n]);\n\n// 4) Responsables -> responsables\nawait replaceParagraphRuns(context, \"Responsables\", [\n  '<w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">Definir </w:t></w:r>',\n  \"<w:r><w:t>responsables</w:t></w:r>\",\n  '<w:r><w:t xml:space=\"preserve\"> de la Calidad</w:t></w:r>'\n]);\n", "ps1": "# \"Deje como estaba el plan de iteracion\"\n# 1) \"Metricas\" -> \"M\u00e9tricas\" (remove spell-check proofErr markers, add accent)\n# 2) \"Herramientas a utilizar\" -> split into \"Describir h\" + \"erramientas a utilizar\"\n#    (visible text becomes \"Describir herramientas a utilizar\")\n# 3) \"Descripcion\" -> \"Descripci\u00f3n\" (remove spell-check proofErr markers, add accent)\n# 4) \"Responsables\" -> \"responsables\" (remove grammar-check proofErr markers, lowercase)\n\n$d = $word.ActiveDocument\n\nfunction Replace-ParagraphRuns {\n    param(\n        [string]$SearchText,\n        [string]$RunsXml\n    )\n    $hitRange = $d.Content\n    $find = $hitRange.Find\n    $find.ClearFormatting()\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($SearchText)\n    if (-not $found) {\n        throw \"Text not found: $SearchText\"\n    }\n    # Rebuild the whole paragraph (not just the matched text) so any sibling\n    # <w:proofErr/> spell/grammar markers inside that paragraph are dropped too.\n    $para = $hitRange.Paragraphs(1).Range\n    $xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' +\n        $RunsXml +\n        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $para.InsertXML($xml)\n}\n\n# 1) Metricas -> M\u00e9tricas\nReplace-ParagraphRuns \"Metricas\" '<w:r><w:t>M\u00e9tricas</w:t></w:r><w:r><w:t xml:space=\"preserve\"> basadas en casos de uso</w:t></w:r>'\n\n# 2) Herramientas a utilizar -> Describir h | erramientas a utilizar\nReplace-ParagraphRuns \"Herramientas a utilizar\" '<w:r><w:t>Describir h</w:t></w:r><w:r><w:t>erramientas a utilizar</w:t></w:r>'\n\n# 3) Descripcion -> Descripci\u00f3n\nReplace-ParagraphRuns \"Descripcion\" '<w:r><w:t>Descripci\u00f3n</w:t></w:r><w:r><w:t xml:space=\"preserve\"> general</w:t></w:r>'\n\n# 4) Responsables -> responsables\nReplace-ParagraphRuns \"Responsables\" '<w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">Definir </w:t></w:r><w:r><w:t>responsables</w:t></w:r><w:r><w:t xml:space=\"preserve\"> de la Calidad</w:t></w:r>'\n"}
